$wb = $excel.ActiveWorkbook

$ws = $wb.Worksheets.Item("ALC")
$ws.Range("H12").Value = 76923440
$ws.Range("I12").Value = 292.85715
$ws.Range("J12").Value = 166667100
$ws.Range("K12").Value = 292.85715
$ws.Range("L12").Value = 166667100
$ws.Range("M12").Value = -122.85715
$ws.Range("N12").Value = -166667440
$ws.Range("H40").Value = 3295.5
$ws.Range("I40").Value = 0
$ws.Range("J40").Value = 3295.5
$ws.Range("K40").Value = 0
$ws.Range("L40").Value = 3295.5
$ws.Range("M40").ClearContents()
$ws.Range("N40").Value = -3645.5
$ws.Range("H101").Value = 860.875
$ws.Range("I101").Value = 898.1429000000001
$ws.Range("J101").Value = 600
$ws.Range("K101").Value = 2694.4287
$ws.Range("L101").Value = 1800
$ws.Range("M101").Value = -1072.4287
$ws.Range("N101").Value = -5044
$ws.Range("H137").Value = 2503.078
$ws.Range("I137").Value = 957.2093
$ws.Range("J137").Value = 3917.383
$ws.Range("K137").Value = 2871.6279
$ws.Range("L137").Value = 11752.149
$ws.Range("M137").Value = -321.6279
$ws.Range("N137").Value = -16852.149
$ws.Range("H138").Value = 1632.79
$ws.Range("I138").Value = 895.55884
$ws.Range("J138").Value = 3199.4062
$ws.Range("K138").Value = 2686.67652
$ws.Range("L138").Value = 9598.2186
$ws.Range("M138").Value = 2453.32348
$ws.Range("N138").Value = -19878.2186
$ws = $wb.Worksheets.Item("ARM")
$ws.Range("H2").Value = 1055.6
$ws.Range("I2").Value = 931.13336
$ws.Range("J2").Value = 1429
$ws.Range("K2").Value = 931.13336
$ws.Range("L2").Value = 1429
$ws.Range("M2").Value = -818.13336
$ws.Range("N2").Value = -1655
$ws.Range("H45").Value = 1237.76
$ws.Range("I45").Value = 819.41174
$ws.Range("J45").Value = 2126.75
$ws.Range("K45").Value = 819.41174
$ws.Range("L45").Value = 2126.75
$ws.Range("M45").Value = -442.41174
$ws.Range("N45").Value = -2880.75
$ws.Range("H74").Value = 3819.4473
$ws.Range("I74").Value = 1074.0968
$ws.Range("J74").Value = 15977.429
$ws.Range("K74").Value = 1074.0968
$ws.Range("L74").Value = 15977.429
$ws.Range("M74").Value = -200.0968
$ws.Range("N74").Value = -17725.429
$ws.Range("H77").Value = 3819.4473
$ws.Range("I77").Value = 1074.0968
$ws.Range("J77").Value = 15977.429
$ws.Range("K77").Value = 5370.484
$ws.Range("L77").Value = 79887.145
$ws.Range("M77").Value = -1002.484
$ws.Range("N77").Value = -88623.145
$ws.Range("H116").Value = 1055.6
$ws.Range("I116").Value = 931.13336
$ws.Range("J116").Value = 1429
$ws.Range("K116").Value = 931.13336
$ws.Range("L116").Value = 1429
$ws.Range("M116").Value = 1362.86664
$ws.Range("N116").Value = -6017
$ws.Range("H132").Value = 5055.7925
$ws.Range("I132").Value = 3597.658
$ws.Range("J132").Value = 8749.733
$ws.Range("K132").Value = 10792.974
$ws.Range("L132").Value = 26249.199
$ws.Range("M132").Value = -8262.974
$ws.Range("N132").Value = -31309.199
$ws = $wb.Worksheets.Item("BSM")
$ws.Range("H3").Value = 1055.6
$ws.Range("I3").Value = 931.13336
$ws.Range("J3").Value = 1429
$ws.Range("K3").Value = 931.13336
$ws.Range("L3").Value = 1429
$ws.Range("M3").Value = -817.13336
$ws.Range("N3").Value = -1657
$ws.Range("H20").Value = 36536.69
$ws.Range("I20").Value = 1514.7646
$ws.Range("J20").Value = 86151.086
$ws.Range("K20").Value = 1514.7646
$ws.Range("L20").Value = 86151.086
$ws.Range("M20").Value = -1267.7646
$ws.Range("N20").Value = -86645.086
$ws.Range("H22").Value = 300.66666
$ws.Range("I22").Value = 216.66667
$ws.Range("K22").Value = 216.66667
$ws.Range("M22").Value = -43.66667000000001
$ws.Range("H129").Value = 43889.5
$ws.Range("J129").Value = 43889.5
$ws.Range("L129").Value = 43889.5
$ws.Range("N129").Value = -53889.5
$ws.Range("H134").Value = 1637.1936
$ws.Range("I134").Value = 1047.5714
$ws.Range("J134").Value = 2875.4
$ws.Range("K134").Value = 3142.7142
$ws.Range("L134").Value = 8626.200000000001
$ws.Range("M134").Value = -607.7142000000003
$ws.Range("N134").Value = -13696.2
$ws = $wb.Worksheets.Item("CRP")
$ws.Range("H31").Value = 11921011
$ws.Range("I31").Value = 35715264
$ws.Range("K31").Value = 35715264
$ws.Range("M31").Value = -35714969
$ws.Range("H34").Value = 11921011
$ws.Range("I34").Value = 35715264
$ws.Range("K34").Value = 35715264
$ws.Range("M34").Value = -35715062
$ws.Range("H58").Value = 22728132
$ws.Range("I58").Value = 41667244
$ws.Range("J58").Value = 1198.9
$ws.Range("K58").Value = 41667244
$ws.Range("L58").Value = 1198.9
$ws.Range("M58").Value = -41667041
$ws.Range("N58").Value = -1604.9
$ws.Range("H94").Value = 1352.8422
$ws.Range("J94").Value = 1528.5714
$ws.Range("L94").Value = 1528.5714
$ws.Range("N94").Value = -2430.5714
$ws.Range("H132").Value = 8824.825999999999
$ws.Range("I132").Value = 12790.9
$ws.Range("J132").Value = 5774
$ws.Range("K132").Value = 38372.7
$ws.Range("L132").Value = 17322
$ws.Range("M132").Value = -35842.7
$ws.Range("N132").Value = -22382
$ws.Range("H136").Value = 22728132
$ws.Range("I136").Value = 41667244
$ws.Range("J136").Value = 1198.9
$ws.Range("K136").Value = 125001732
$ws.Range("L136").Value = 3596.7
$ws.Range("M136").Value = -124999182
$ws.Range("N136").Value = -8696.700000000001
$ws = $wb.Worksheets.Item("CUL")
$ws.Range("H122").Value = 887.30554
$ws.Range("I122").Value = 405.57144
$ws.Range("J122").Value = 1561.7333
$ws.Range("K122").Value = 3650.14296
$ws.Range("L122").Value = 14055.5997
$ws.Range("M122").Value = -1200.14296
$ws.Range("N122").Value = -18955.5997
$ws.Range("H132").Value = 52185.684
$ws.Range("J132").Value = 111710.5
$ws.Range("L132").Value = 1005394.5
$ws.Range("N132").Value = -1010454.5
$ws = $wb.Worksheets.Item("GSM")
$ws.Range("H48").Value = 14850
$ws.Range("J48").Value = 14850
$ws.Range("L48").Value = 14850
$ws.Range("N48").Value = -15820
$ws.Range("H70").Value = 5000
$ws.Range("I70").Value = 5000
$ws.Range("J70").Value = 5000
$ws.Range("K70").Value = 5000
$ws.Range("L70").Value = 5000
$ws.Range("M70").Value = -4730
$ws.Range("N70").Value = -5540
$ws.Range("H73").Value = 5000
$ws.Range("I73").Value = 5000
$ws.Range("J73").Value = 5000
$ws.Range("K73").Value = 5000
$ws.Range("L73").Value = 5000
$ws.Range("M73").Value = -4064
$ws.Range("N73").Value = -6872
$ws.Range("H96").Value = 13130.5
$ws.Range("J96").Value = 13130.5
$ws.Range("L96").Value = 13130.5
$ws.Range("N96").Value = -18622.5
$ws = $wb.Worksheets.Item("LTW")
$ws.Range("H7").Value = 1363.4667
$ws.Range("I7").Value = 1078.5
$ws.Range("J7").Value = 2503.3333
$ws.Range("K7").Value = 1078.5
$ws.Range("L7").Value = 2503.3333
$ws.Range("M7").Value = -966.5
$ws.Range("N7").Value = -2727.3333
$ws.Range("H16").Value = 632
$ws.Range("I16").Value = 0
$ws.Range("J16").Value = 632
$ws.Range("K16").Value = 0
$ws.Range("L16").Value = 632
$ws.Range("M16").ClearContents()
$ws.Range("N16").Value = -972
$ws.Range("H46").Value = 1469.1
$ws.Range("I46").Value = 1265
$ws.Range("K46").Value = 1265
$ws.Range("M46").Value = -1077
$ws.Range("H126").Value = 1363.4667
$ws.Range("I126").Value = 1078.5
$ws.Range("J126").Value = 2503.3333
$ws.Range("K126").Value = 3235.5
$ws.Range("L126").Value = 7509.999899999999
$ws.Range("M126").Value = -765.5
$ws.Range("N126").Value = -12449.9999
$ws.Range("H129").Value = 37429
$ws.Range("J129").Value = 37429
$ws.Range("L129").Value = 37429
$ws.Range("N129").Value = -47429
$ws.Range("H132").Value = 5862.6772
$ws.Range("I132").Value = 7609.4736
$ws.Range("J132").Value = 3096.9167
$ws.Range("K132").Value = 22828.4208
$ws.Range("L132").Value = 9290.750100000001
$ws.Range("M132").Value = -20298.4208
$ws.Range("N132").Value = -14350.7501

Write-Host "applied 208 changes"